$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data: cell A3 gets the new text value.
$ws.Range("A3").Value = "добавлена строка"

# Move the selection to A4, as if the user pressed Enter after typing into A3.
$ws.Range("A4").Select()
